$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.970.02'
Set-TextValue $ws.Range('E2') '  -0.99%  '
Set-TextValue $ws.Range('D3') '1.810.86'
Set-TextValue $ws.Range('E3') '  -0.49%  '
Set-TextValue $ws.Range('E4') '  +0.05%  '
Set-TextValue $ws.Range('D5') '310.89'
Set-TextValue $ws.Range('E5') '  -0.68%  '
Set-TextValue $ws.Range('D6') '1.002'
Set-TextValue $ws.Range('E6') '  +0.10%  '
Set-TextValue $ws.Range('D7') '0.4632'
Set-TextValue $ws.Range('E7') '  +4.09%  '
Set-TextValue $ws.Range('D8') '0.3720'
Set-TextValue $ws.Range('E8') '  -1.11%  '
Set-TextValue $ws.Range('D9') '0.07378'
Set-TextValue $ws.Range('E9') '  -0.13%  '
Set-TextValue $ws.Range('D10') '0.8757'
Set-TextValue $ws.Range('E10') '  -0.35%  '
Set-TextValue $ws.Range('D11') '20.47'
Set-TextValue $ws.Range('E11') '  -1.68%  '
Set-TextValue $ws.Range('D12') '1.818.81'
Set-TextValue $ws.Range('E12') '  -0.04%  '
Set-TextValue $ws.Range('D13') '5.372'
Set-TextValue $ws.Range('E13') '  -0.81%  '
Set-TextValue $ws.Range('D14') '92.52'
Set-TextValue $ws.Range('E14') '  -0.58%  '
Set-TextValue $ws.Range('D15') '6.532'
Set-TextValue $ws.Range('E15') '  -2.46%  '
Set-TextValue $ws.Range('D16') '0.07053'
Set-TextValue $ws.Range('E16') '  -0.75%  '
Set-TextValue $ws.Range('E17') '  +0.05%  '
Set-TextValue $ws.Range('D18') '0.000008720'
Set-TextValue $ws.Range('E18') '  -0.87%  '
Set-TextValue $ws.Range('E19') '  +0.17%  '
Set-TextValue $ws.Range('D20') '14.72'
Set-TextValue $ws.Range('E20') '  -1.90%  '
Set-TextValue $ws.Range('D21') '26.981.22'
Set-TextValue $ws.Range('E21') '  -1.00%  '
Set-TextValue $ws.Range('D22') '5.317'
Set-TextValue $ws.Range('E22') '  -0.59%  '
Set-TextValue $ws.Range('D23') '10.65'
Set-TextValue $ws.Range('E23') '  -2.38%  '
Set-TextValue $ws.Range('D24') '2.068.88'
Set-TextValue $ws.Range('E24') '  +1.07%  '
Set-TextValue $ws.Range('D25') '1.903'
Set-TextValue $ws.Range('E25') '  -2.79%  '
Set-TextValue $ws.Range('E26') '  +0.43%  '
Set-TextValue $ws.Range('D27') '18.40'
Set-TextValue $ws.Range('E27') '  -0.81%  '
Set-TextValue $ws.Range('D28') '2.147'
Set-TextValue $ws.Range('E28') '  -6.45%  '
Set-TextValue $ws.Range('D29') '5.284'
Set-TextValue $ws.Range('E29') '  -1.03%  '
Set-TextValue $ws.Range('D30') '116.00'
Set-TextValue $ws.Range('E30') '  -1.00%  '
Set-TextValue $ws.Range('D31') '0.08942'
Set-TextValue $ws.Range('E31') '  +0.90%  '
Set-TextValue $ws.Range('D32') '0.7572'
Set-TextValue $ws.Range('E32') '  -3.69%  '
Set-TextValue $ws.Range('D33') '1.158'
Set-TextValue $ws.Range('E33') '  -2.85%  '
Set-TextValue $ws.Range('D34') '2.930'
Set-TextValue $ws.Range('E34') '  +0.58%  '
Set-TextValue $ws.Range('E35') '  -2.08%  '
Set-TextValue $ws.Range('E36') '  +0.10%  '
Set-TextValue $ws.Range('D37') '1.105'
Set-TextValue $ws.Range('D38') '0.01979'
Set-TextValue $ws.Range('E38') '  +0.67%  '
Set-TextValue $ws.Range('B39') 'RenderToken'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D39') '2.452'
Set-TextValue $ws.Range('E39') '  +7.47%  '
Set-TextValue $ws.Range('B40') 'Hedera'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D40') '0.05245'
Set-TextValue $ws.Range('E40') '  -0.25%  '
Set-TextValue $ws.Range('D41') '2.921'
Set-TextValue $ws.Range('E41') '  +1.96%  '
Set-TextValue $ws.Range('B42') 'TheSandbox'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D42') '0.5310'
Set-TextValue $ws.Range('E42') '  +0.57%  '
Set-TextValue $ws.Range('B43') 'FraxShare'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '7.201'
Set-TextValue $ws.Range('E43') '  -1.30%  '
Set-TextValue $ws.Range('D44') '0.1666'
Set-TextValue $ws.Range('E44') '  -2.01%  '
Set-TextValue $ws.Range('D45') '8.526'
Set-TextValue $ws.Range('E45') '  -0.77%  '
Set-TextValue $ws.Range('D46') '0.4990'
Set-TextValue $ws.Range('E46') '  -0.72%  '
Set-TextValue $ws.Range('D47') '10.39'
Set-TextValue $ws.Range('E47') '  -2.08%  '
Set-TextValue $ws.Range('E48') '  -0.50%  '
Set-TextValue $ws.Range('B49') 'Quant'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D49') '104.26'
Set-TextValue $ws.Range('E49') '  -0.60%  '
Set-TextValue $ws.Range('B50') 'PaxDollar'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D50') '1.002'
Set-TextValue $ws.Range('E50') '  +0.14%  '
Set-TextValue $ws.Range('D51') '0.06298'
Set-TextValue $ws.Range('E51') '  -1.28%  '
